{"js": "// 1) Update the cached DATE field result from \"15 September 2023\" to\n//    \"22 September 2023\". The field's cached text lives in a plain run\n//    (<w:t>15 September 2023</w:t>), so a literal text search/replace on\n//    the document body targets it precisely.\nconst body = context.document.body;\nconst dateResults = body.search(\"15 September 2023\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"22 September 2023\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Append a new bullet/list item right after the paragraph that ends\n//    with \"...network connectivity is restored.\" and before the\n//    \"CONCLUSION\" heading. The new paragraph continues the same\n//    numbered list (ListParagraph style / numId 10) as its predecessor.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"network connectivity is restored\") !== -1) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (anchorParagraph) {\n  const newParagraph = anchorParagraph.insertParagraph(\n    \"Hiccups in network connectivity were discovered to cause the streamers to briefly lose connection before \",\n    \"After\"\n  );\n  newParagraph.insertText(\"quickly coming back online. \", \"End\");\n  newParagraph.insertText(\"The monitor was \", \"End\");\n  newParagraph.insertText(\"modified to \", \"End\");\n  newParagraph.insertText(\"double \", \"End\");\n  newParagraph.insertText(\"check \", \"End\");\n  newParagraph.insertText(\"if the streamers were offline for a period of \", \"End\");\n  newParagraph.insertText(\"10 \", \"End\");\n  newParagraph.insertText(\"seconds\", \"End\");\n  newParagraph.insertText(\" before sending a notification.\", \"End\");\n  await context.sync();\n}\n", "ps1": "# 1) Update the cached DATE field result from \"15 September 2023\" to\n#    \"22 September 2023\". The field's cached text lives in a plain run\n#    (<w:t>15 September 2023</w:t>), so a literal Find/Replace on the\n#    document content targets it precisely.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"15 September 2023\",  # FindText\n    $false,                 # MatchCase\n    $false,                 # MatchWholeWord\n    $false,                 # MatchWildcards\n    $false,                 # MatchSoundsLike\n    $false,                 # MatchAllWordForms\n    $true,                  # Forward\n    1,                       # Wrap (wdFindContinue)\n    $false,                 # Format\n    \"22 September 2023\",    # ReplaceWith\n    2                        # Replace (wdReplaceOne)\n) | Out-Null\n\n# 2) Append a new bullet/list item right after the paragraph that ends\n#    with \"...network connectivity is restored.\" and before the\n#    \"CONCLUSION\" heading. The new paragraph continues the same\n#    numbered list (ListParagraph style / numId 10) as its predecessor.\n$anchorIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*network connectivity is restored*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -gt 0) {\n    $d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter()\n\n    $newParaIndex = $anchorIndex + 1\n    $parts = @(\n        \"Hiccups in network connectivity were discovered to cause the streamers to briefly lose connection before \",\n        \"quickly coming back online. \",\n        \"The monitor was \",\n        \"modified to \",\n        \"double \",\n        \"check \",\n        \"if the streamers were offline for a period of \",\n        \"10 \",\n        \"seconds\",\n        \" before sending a notification.\"\n    )\n    foreach ($part in $parts) {\n        # Re-fetch the paragraph's Range each time so the insertion point is\n        # re-anchored at the (now longer) paragraph's true end, immediately\n        # before its paragraph mark.\n        $d.Paragraphs.Item($newParaIndex).Range.InsertAfter($part)\n    }\n}\n"}
